# Adds a new "2022-Q1" sheet (with per-fund holding detail) right before the
# "总计" (totals) sheet, and prepends a matching summary row to "总计".

function Set-TextCell($cell, $value) {
    # Force the cell to be stored as text (keeps leading zeros / numeric-looking
    # strings like "180.12" as text) while resetting the style back to the
    # workbook default so no stray style index gets attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-FormatFrom($srcCell, $dstCell) {
    # Copies only the formatting (style) of $srcCell onto $dstCell.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before the last sheet ("总计").
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Reference sheet that already has the per-fund layout/styling we want to copy
# (any of the 2021-* detail sheets works; use the last one, "2021-Q4").
$refSheet = $wb.Worksheets.Item("2021-Q4")

# Header row (B1:H1), bold+bordered style copied from the reference sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $dst = $newSheet.Cells.Item(1, $col)
    Set-FormatFrom $refSheet.Cells.Item(1, $col) $dst
    $dst.Value = $headers[$col - 2]
}

$fundData = @(
    @("001475", "易方达国防军工混合", "180.12", "92.18", "4.63", "8.3396", 8),
    @("004139", "中邮军民融合灵活配置混合", "16.83", "86.35", "8.87", "1.4928", 1),
    @("002983", "长信国防军工量化灵活配置混合A", "26.74", "93.30", "5.38", "1.4386", 7),
    @("110005", "易方达积极成长混合", "36.95", "92.57", "3.44", "1.2711", 10),
    @("001224", "中邮新思路灵活配置混合", "33.47", "74.78", "3.20", "1.0710", 9),
    @("040015", "华安动态灵活配置混合", "22.93", "79.55", "4.57", "1.0479", 3),
    @("010792", "华安成长先锋混合A", "16.44", "93.43", "4.67", "0.7677", 7),
    @("008960", "长信国防军工量化灵活配置混合C", "13.19", "93.30", "5.38", "0.7096", 7),
    @("006154", "华安制造先锋混合", "14.05", "93.81", "4.42", "0.6210", 7),
    @("001479", "中邮风格轮动灵活配置混合", "9.45", "62.17", "5.22", "0.4933", 1),
    @("008980", "中邮科技创新精选混合A", "9.42", "76.05", "4.39", "0.4135", 3),
    @("580009", "东吴多策略灵活配置混合", "8.31", "88.46", "4.53", "0.3764", 8),
    @("160425", "华安创业板两年定期开放混合", "5.11", "96.75", "4.84", "0.2473", 1),
    @("008009", "华商高端装备制造股票", "5.40", "88.04", "3.76", "0.2030", 6),
    @("005457", "景顺长城量化小盘股票", "9.49", "93.39", "1.96", "0.1860", 4),
    @("010793", "华安成长先锋混合C", "3.80", "93.43", "4.67", "0.1775", 7),
    @("008981", "中邮科技创新精选混合C", "4.02", "76.05", "4.39", "0.1765", 3),
    @("001047", "光大保德信国企改革主题股票", "2.87", "85.82", "4.85", "0.1392", 5),
    @("310368", "申万菱信竞争优势混合", "0.83", "91.22", "4.45", "0.0369", 4),
    @("163818", "中银中小盘成长混合", "0.98", "87.49", "2.18", "0.0214", 9),
    @("008082", "国寿安保研究精选混合A", "0.52", "91.60", "3.56", "0.0185", 10),
    @("519971", "长信改革红利灵活配置混合", "0.27", "73.88", "3.23", "0.0087", 9),
    @("008083", "国寿安保研究精选混合C", "0.15", "91.60", "3.56", "0.0053", 10),
    @("519165", "新华鑫利灵活配置混合", "0.05", "74.84", "2.92", "0.0015", 9)
)

$row = 2
foreach ($rec in $fundData) {
    $aCell = $newSheet.Cells.Item($row, 1)
    Set-FormatFrom $refSheet.Cells.Item($row, 1) $aCell
    $aCell.Value = $row - 2

    Set-TextCell $newSheet.Cells.Item($row, 2) $rec[0]
    Set-TextCell $newSheet.Cells.Item($row, 3) $rec[1]
    Set-TextCell $newSheet.Cells.Item($row, 4) $rec[2]
    Set-TextCell $newSheet.Cells.Item($row, 5) $rec[3]
    Set-TextCell $newSheet.Cells.Item($row, 6) $rec[4]
    Set-TextCell $newSheet.Cells.Item($row, 7) $rec[5]
    $newSheet.Cells.Item($row, 8).Value = $rec[6]

    $row++
}

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).Style = "Normal"

Set-FormatFrom $totalSheet.Cells.Item(3, 1) $totalSheet.Cells.Item(2, 1)
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 24
$totalSheet.Cells.Item(2, 4).Value = 19.26

# Renumber the index column (A) for the rows that shifted down by one.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
